$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("C12").Value = 2.55
$ws.Range("D12").Value = 2.083333333333333
$ws.Range("C13").Value = 3.261538461538461
$ws.Range("D13").Value = 2.846153846153846
$ws.Range("C14").Value = 5.32
$ws.Range("D14").Value = 5.4
$ws.Range("C15").Value = 6.516666666666667
$ws.Range("D15").Value = 6.083333333333333
$ws.Range("C16").Value = 7.18
$ws.Range("D16").Value = 7.3
$ws.Range("C17").Value = 6.659999999999999
$ws.Range("D17").Value = 6.1
$ws.Range("C18").Value = 7.230769230769231
$ws.Range("D18").Value = 8.076923076923077
$ws.Range("C19").Value = 6.985714285714286
$ws.Range("D19").Value = 8
$ws.Range("C20").Value = 7.583333333333333
$ws.Range("D20").Value = 7.75
$ws.Range("C21").Value = 7.533333333333333
$ws.Range("D21").Value = 7.333333333333333
$ws.Range("C22").Value = 7.111111111111111
$ws.Range("D22").Value = 6.777777777777778
$ws.Range("C23").Value = 5.0625
$ws.Range("D23").Value = 4.6875
$ws.Range("C24").Value = 1.181818181818182
$ws.Range("D24").Value = 1.636363636363636
$ws.Range("D35").Value = 0.1111111111111111
$ws.Range("C36").Value = 2.036363636363636
$ws.Range("D36").Value = 2.090909090909091
$ws.Range("C37").Value = 3.447058823529412
$ws.Range("D37").Value = 3.058823529411764
$ws.Range("C38").Value = 5.74
$ws.Range("D38").Value = 5.3
$ws.Range("C39").Value = 5.876923076923076
$ws.Range("D39").Value = 5.615384615384615
$ws.Range("C40").Value = 5.8
$ws.Range("D40").Value = 5.875
$ws.Range("C41").Value = 7
$ws.Range("D41").Value = 6.545454545454546
$ws.Range("C42").Value = 6.675
$ws.Range("D42").Value = 8.375
$ws.Range("C43").Value = 6.893333333333334
$ws.Range("D43").Value = 7.133333333333334
$ws.Range("C44").Value = 7.188888888888889
$ws.Range("D44").Value = 7
$ws.Range("C45").Value = 7.44
$ws.Range("D45").Value = 7
$ws.Range("C46").Value = 7
$ws.Range("D46").Value = 6.533333333333333
$ws.Range("C47").Value = 5.127272727272727
$ws.Range("D47").Value = 4.727272727272728
$ws.Range("C48").Value = 1.054545454545454
$ws.Range("D48").Value = 1.636363636363636
$ws.Range("C49").Value = 0.1333333333333333
$ws.Range("C60").Value = 2.285714285714286
$ws.Range("D60").Value = 2.357142857142857
$ws.Range("C61").Value = 3.163636363636364
$ws.Range("D61").Value = 3
$ws.Range("C62").Value = 5.644444444444444
$ws.Range("D62").Value = 5
$ws.Range("C63").Value = 6.214285714285714
$ws.Range("D63").Value = 6.214285714285714
$ws.Range("C64").Value = 6.333333333333333
$ws.Range("D64").Value = 6.444444444444445
$ws.Range("C65").Value = 6.428571428571429
$ws.Range("D65").Value = 6.071428571428571
$ws.Range("C66").Value = 6.957142857142857
$ws.Range("D66").Value = 8.214285714285714
$ws.Range("C67").Value = 7.2875
$ws.Range("D67").Value = 7.75
$ws.Range("C68").Value = 7.092307692307693
$ws.Range("D68").Value = 7.153846153846154
$ws.Range("C69").Value = 6.733333333333333
$ws.Range("D69").Value = 6.5
$ws.Range("C70").Value = 6.183333333333334
$ws.Range("D70").Value = 5.583333333333333
$ws.Range("C71").Value = 4.54
$ws.Range("D71").Value = 4.3
$ws.Range("C72").Value = 1.414285714285714
$ws.Range("D72").Value = 1.428571428571429
$ws.Range("C73").Value = 0.1764705882352941
$ws.Range("D73").Value = 0.05882352941176471
$ws.Range("C84").Value = 2.636363636363636
$ws.Range("D84").Value = 2.272727272727273
$ws.Range("C85").Value = 3.566666666666666
$ws.Range("D85").Value = 3.083333333333333
$ws.Range("C86").Value = 5.24
$ws.Range("D86").Value = 4.6
$ws.Range("C87").Value = 5.966666666666666
$ws.Range("D87").Value = 5.916666666666667
$ws.Range("C88").Value = 7.05
$ws.Range("D88").Value = 7.25
$ws.Range("C89").Value = 6.984615384615385
$ws.Range("D89").Value = 6.538461538461538
$ws.Range("C90").Value = 7.62
$ws.Range("D90").Value = 8.199999999999999
$ws.Range("C91").Value = 8.015384615384615
$ws.Range("D91").Value = 8.692307692307692
$ws.Range("C92").Value = 8.18
$ws.Range("C93").Value = 8.18
$ws.Range("D93").Value = 7.7
$ws.Range("C94").Value = 6.716666666666666
$ws.Range("D94").Value = 6.916666666666667
$ws.Range("C95").Value = 6.12
$ws.Range("D95").Value = 5.6
$ws.Range("C96").Value = 1.45
$ws.Range("D96").Value = 1.916666666666667
$ws.Range("C97").Value = 0.08888888888888889
$ws.Range("C108").Value = 3.309090909090909
$ws.Range("D108").Value = 3.363636363636364
$ws.Range("C109").Value = 4.64
$ws.Range("D109").Value = 4.5
$ws.Range("C110").Value = 6.466666666666667
$ws.Range("D110").Value = 6.777777777777778
$ws.Range("C111").Value = 7.127272727272728
$ws.Range("D111").Value = 7.454545454545454
$ws.Range("C112").Value = 8.036363636363637
$ws.Range("D112").Value = 7.363636363636363
$ws.Range("C113").Value = 8.254545454545454
$ws.Range("D113").Value = 7.272727272727272
$ws.Range("C114").Value = 8.82
$ws.Range("D114").Value = 9.800000000000001
$ws.Range("C115").Value = 8.4
$ws.Range("D115").Value = 8.818181818181818
$ws.Range("C116").Value = 8.199999999999999
$ws.Range("D116").Value = 8.066666666666666
$ws.Range("C117").Value = 6.8
$ws.Range("D117").Value = 6.214285714285714
$ws.Range("C118").Value = 6.107692307692307
$ws.Range("D118").Value = 5.846153846153846
$ws.Range("C119").Value = 3.844444444444445
$ws.Range("D119").Value = 4
$ws.Range("C120").Value = 1.78
$ws.Range("D120").Value = 1.4
$ws.Range("C121").Value = 0.16
$ws.Range("D131").Value = 0
$ws.Range("C132").Value = 2.171428571428572
$ws.Range("D132").Value = 1.857142857142857
$ws.Range("C133").Value = 2.4
$ws.Range("D133").Value = 2.142857142857143
$ws.Range("C134").Value = 4.711111111111111
$ws.Range("D134").Value = 4.666666666666667
$ws.Range("C135").Value = 5.4
$ws.Range("D135").Value = 4.8
$ws.Range("C136").Value = 5.885714285714286
$ws.Range("D136").Value = 5.857142857142857
$ws.Range("C137").Value = 5.866666666666666
$ws.Range("D137").Value = 5.933333333333334
$ws.Range("C138").Value = 6.646153846153847
$ws.Range("D138").Value = 7.846153846153846
$ws.Range("C139").Value = 6.463157894736842
$ws.Range("D139").Value = 7.052631578947368
$ws.Range("C140").Value = 5.76
$ws.Range("D140").Value = 6.2
$ws.Range("C141").Value = 6.727272727272728
$ws.Range("D141").Value = 6.727272727272728
$ws.Range("C142").Value = 5.685714285714285
$ws.Range("D142").Value = 5.142857142857143
$ws.Range("C143").Value = 4.46
$ws.Range("D143").Value = 4.2
$ws.Range("C144").Value = 1.104761904761905
$ws.Range("D144").Value = 1.380952380952381
$ws.Range("C145").Value = 0.1818181818181818
$ws.Range("C155").Value = 0.06
$ws.Range("D155").Value = 0.2
$ws.Range("C156").Value = 3.725
$ws.Range("D156").Value = 3.375
$ws.Range("C157").Value = 5.171428571428572
$ws.Range("D157").Value = 5.142857142857143
$ws.Range("C158").Value = 6.58
$ws.Range("D158").Value = 6.6
$ws.Range("C159").Value = 7.857142857142857
$ws.Range("D159").Value = 7.428571428571429
$ws.Range("C160").Value = 8.012499999999999
$ws.Range("D160").Value = 7.8125
$ws.Range("C161").Value = 9.054545454545455
$ws.Range("D161").Value = 8.454545454545455
$ws.Range("C162").Value = 9.461538461538462
$ws.Range("D162").Value = 11.15384615384615
$ws.Range("C163").Value = 8.854545454545455
$ws.Range("D163").Value = 9.363636363636363
$ws.Range("C164").Value = 9.26
$ws.Range("D164").Value = 9.699999999999999
$ws.Range("C165").Value = 8.6
$ws.Range("D165").Value = 8.461538461538462
$ws.Range("C166").Value = 7.466666666666667
$ws.Range("C167").Value = 5.866666666666667
$ws.Range("D167").Value = 5.416666666666667
$ws.Range("C168").Value = 1.575
$ws.Range("D168").Value = 1.75
$ws.Range("C169").Value = 0.1733333333333333
